# "Added Indian MF 1st Stab" - weekly MarketBeat rank refresh.
# Each week a new column is inserted right after column A (i.e. at B),
# pushing all the previously-collected weekly columns one slot to the
# right.  This commit brings the sheet up to date by nine weeks in one
# shot (Jun_16, Jun_24, Jun_30, Jul_07, Jul_17, Jul_23, Aug_04, Aug_25,
# Sep_08), defaulting the freshly inserted cells to "UN" (unchanged)
# except where an actual rating action happened that week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert nine new blank columns at B, shifting the old B:V block
#    (Jun_09 .. Mar_10) out to K:AE.
# ---------------------------------------------------------------
for ($i = 0; $i -lt 9; $i++) {
    $ws.Columns("B").Insert()
}

# ---------------------------------------------------------------
# 2. New week-ending header dates for the freshly inserted columns.
# ---------------------------------------------------------------
$newDates = @{
    "B1" = "Sep_08"
    "C1" = "Aug_25"
    "D1" = "Aug_04"
    "E1" = "Jul_23"
    "F1" = "Jul_17"
    "G1" = "Jul_07"
    "H1" = "Jun_30"
    "I1" = "Jun_24"
    "J1" = "Jun_16"
}
foreach ($addr in $newDates.Keys) {
    $ws.Range($addr).Value = $newDates[$addr]
}

# ---------------------------------------------------------------
# 3. Default every data row's nine new cells (columns B:J) to "UN"
#    (no analyst activity that week) for every row that already had
#    data before this edit (rows 2 through 33).
# ---------------------------------------------------------------
$lastRow = 33
for ($r = 2; $r -le $lastRow; $r++) {
    $rng = $ws.Range("B" + $r + ":J" + $r)
    $rng.Value = "UN"
}

# ---------------------------------------------------------------
# 4. Real rating-change events that landed in the new weeks.
# ---------------------------------------------------------------
$ws.Range("C7").Value  = "8/23/2019,Lowers Target,Equal Weight,$34.00 -> $18.00"
$ws.Range("C18").Value = "8/23/2019,Downgrades,Outperform -> Sector Perform,$30.00 -> $22.00"
$ws.Range("C19").Value = "8/23/2019,Lowers Target,Neutral,$27.00 -> $21.00"
$ws.Range("F19").Value = "7/15/2019,Downgrades,Buy -> Neutral,$26.50 -> $27.00"
$ws.Range("J27").Value = "6/11/2019,Downgrades,Overweight -> Equal Weight,C$45.00 -> C$46.00"

# Highlight the cells carrying this week's new rating-change events
# (same convention the sheet already uses for the oldest column, now
# shifted to AE, where V10/V18/V20/V21 were tinted green).
$highlightAddrs = @("C7", "C18", "C19", "F19", "J27")
foreach ($addr in $highlightAddrs) {
    $ws.Range($addr).Interior.Color = 13408767
}
